$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / country-name updates (shared string reordering effects) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 01:22"
$ws.Range("A27").Value = "Chequia"
$ws.Range("A28").Value = "Chile"
$ws.Range("A29").Value = "India"
$ws.Range("A133").Value = "Guayana Francesa"
$ws.Range("A134").Value = "Aruba"
$ws.Range("A135").Value = "Guatemala"
$ws.Range("A136").Value = "El Salvador"

# --- Numeric value updates ---
$ws.Range("B4").Value = 364088
$ws.Range("C4").Value = 27415
$ws.Range("E4").Value = 333760
$ws.Range("B24").Value = 5895
$ws.Range("C24").Value = 145
$ws.Range("E24").Value = 3422
$ws.Range("B25").Value = 5865
$ws.Range("C25").Value = 178
$ws.Range("E25").Value = 5757
$ws.Range("B27").Value = 4822
$ws.Range("C27").Value = 235
$ws.Range("D27").Value = 121
$ws.Range("E27").Value = 4623
$ws.Range("F27").Value = 84
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = 78
$ws.Range("B28").Value = 4815
$ws.Range("C28").Value = 344
$ws.Range("D28").Value = 728
$ws.Range("E28").Value = 4050
$ws.Range("F28").Value = 327
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 37
$ws.Range("B29").Value = 4778
$ws.Range("C29").Value = 489
$ws.Range("D29").Value = 375
$ws.Range("E29").Value = 4267
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 18
$ws.Range("H29").Value = 136
$ws.Range("B133").Value = 72
$ws.Range("C133").Value = 4
$ws.Range("D133").Value = 34
$ws.Range("E133").Value = 38
$ws.Range("F133").Value = 1
$ws.Range("B134").Value = 71
$ws.Range("C134").Value = 7
$ws.Range("D134").Value = 2
$ws.Range("E134").Value = 69
$ws.Range("F134").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("B135").Value = 70
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 15
$ws.Range("E135").Value = 52
$ws.Range("F135").Value = 3
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 3
$ws.Range("B136").Value = 69
$ws.Range("C136").Value = 7
$ws.Range("D136").Value = 5
$ws.Range("E136").Value = 60
$ws.Range("F136").Value = 4
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 4
$ws.Range("F144").Value = 1
